# added harvard case classification
# Fill in the previously-blank "_old" model columns (Ada_old=C, Avey_old=F,
# K health_old=M, WebMD_old=Q, doctor_MA_old=S, doctor_NJ_old=U,
# doctor_TH_old=W) with their classification metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: precision
$ws.Range("C2").Value = 0.6666666666666666
$ws.Range("F2").Value = 1
$ws.Range("M2").Value = 0.5
$ws.Range("Q2").Value = 0.3333333333333333
$ws.Range("S2").Value = 0.6666666666666666
$ws.Range("U2").Value = 1
$ws.Range("W2").Value = 1

# Row 3: recall
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("M3").Value = 0.5
$ws.Range("Q3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("U3").Value = 0.5
$ws.Range("W3").Value = 1

# Row 4: f1-score
$ws.Range("C4").Value = 0.8
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.5
$ws.Range("Q4").Value = 0.5
$ws.Range("S4").Value = 0.8
$ws.Range("U4").Value = 0.6666666666666666
$ws.Range("W4").Value = 1

# Row 5: f2-score
$ws.Range("C5").Value = 0.9090909090909091
$ws.Range("F5").Value = 0.5555555555555556
$ws.Range("M5").Value = 0.5
$ws.Range("Q5").Value = 0.7142857142857143
$ws.Range("S5").Value = 0.9090909090909091
$ws.Range("U5").Value = 0.5555555555555556
$ws.Range("W5").Value = 1

# Row 6: NDCG
$ws.Range("C6").Value = 1
$ws.Range("F6").Value = 0.8262346571285599
$ws.Range("M6").Value = 0.8262346571285599
$ws.Range("Q6").Value = 0.5950427489208391
$ws.Range("S6").Value = 0.9639404333166532
$ws.Range("U6").Value = 0.8262346571285599
$ws.Range("W6").Value = 1

# Row 7: M1 (boolean)
$ws.Range("C7").Value = $true
$ws.Range("F7").Value = $true
$ws.Range("M7").Value = $true
$ws.Range("S7").Value = $true
$ws.Range("U7").Value = $true
$ws.Range("W7").Value = $true

# Row 8: M3 (boolean)
$ws.Range("C8").Value = $true
$ws.Range("F8").Value = $true
$ws.Range("M8").Value = $true
$ws.Range("S8").Value = $true
$ws.Range("U8").Value = $true
$ws.Range("W8").Value = $true

# Row 9: M5 (boolean)
$ws.Range("C9").Value = $true
$ws.Range("F9").Value = $true
$ws.Range("M9").Value = $true
$ws.Range("Q9").Value = $true
$ws.Range("S9").Value = $true
$ws.Range("U9").Value = $true
$ws.Range("W9").Value = $true

# Row 10: position
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("M10").Value = 1
$ws.Range("Q10").Value = 5
$ws.Range("S10").Value = 1
$ws.Range("U10").Value = 1
$ws.Range("W10").Value = 1
